# store.xlsx update: add a second "store" (admin id / value) column.
#
# Sheet1 originally has 4 columns (A:D):
#   A=TestCase, B=store, C=store, D=wait
# A new "store" column is inserted before the old column D (old D -> E),
# duplicating column C's header/format and adding a new admin2 JSON row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D; existing D ("wait") shifts right to E.
$ws.Columns.Item(4).Insert()

# Match the new column's width to column C's width (both are "store" columns).
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(3).ColumnWidth

# Copy column C's formatting (header style, borders, fonts, etc.) into the
# new column D, row by row, then fill in the new column's values.

# Row 1 header: same "store" header/style as column C.
$ws.Cells.Item(1, 3).Copy() | Out-Null
$ws.Cells.Item(1, 4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1, 4).Value = "store"

# Row 2: left blank, matching column C/B's data-row style.
$ws.Cells.Item(2, 3).Copy() | Out-Null
$ws.Cells.Item(2, 4).PasteSpecial(-4122) | Out-Null

# Row 3: new admin2 JSON payload, matching column C/B's data-row style.
$ws.Cells.Item(3, 3).Copy() | Out-Null
$ws.Cells.Item(3, 4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(3, 4).Value = '{"target":"admin2","value":"name2"}'

$excel.CutCopyMode = $false

# Restore the selection that was active when the workbook was last saved.
$ws.Range("D5").Select() | Out-Null
